# Update LR-pair edge statistics (Il10-Il10ra) with recomputed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E:T for data rows 2..17, built from the new TPM-based NATMI output.
$newValues = New-Object 'object[,]' 16,16

# Row 2
$newValues[0,0] = 3
$newValues[0,1] = 1
$newValues[0,2] = 1.208813
$newValues[0,3] = 3.626439
$newValues[0,4] = 0.06175614228936263
$newValues[0,5] = 0.06175614228936264
$newValues[0,6] = 3
$newValues[0,7] = 1
$newValues[0,8] = 0.1956483333333333
$newValues[0,9] = 0.5869450000000001
$newValues[0,10] = 0.003542563479678112
$newValues[0,11] = 0.003542563479678113
$newValues[0,12] = 0.2365022487616667
$newValues[0,13] = 2.128520238855001
$newValues[0,14] = 0.0002187750543201011
$newValues[0,15] = 0.0002187750543201011

# Row 3
$newValues[1,0] = 3
$newValues[1,1] = 1
$newValues[1,2] = 1.208813
$newValues[1,3] = 3.626439
$newValues[1,4] = 0.06175614228936263
$newValues[1,5] = 0.06175614228936264
$newValues[1,6] = 1
$newValues[1,7] = 0.3333333333333333
$newValues[1,8] = 0.07822766666666667
$newValues[1,9] = 0.234683
$newValues[1,10] = 0.001416452010156485
$newValues[1,11] = 0.001416452010156486
$newValues[1,12] = 0.09456262042633334
$newValues[1,13] = 0.8510635838370001
$newValues[1,14] = 0.00008747461188527763
$newValues[1,15] = 0.00008747461188527766

# Row 4
$newValues[2,0] = 3
$newValues[2,1] = 1
$newValues[2,2] = 1.208813
$newValues[2,3] = 3.626439
$newValues[2,4] = 0.06175614228936263
$newValues[2,5] = 0.06175614228936264
$newValues[2,6] = 2
$newValues[2,7] = 0.6666666666666666
$newValues[2,8] = 0.04544366666666667
$newValues[2,9] = 0.136331
$newValues[2,10] = 0.0008228389742616373
$newValues[2,11] = 0.0008228389742616374
$newValues[2,12] = 0.05493289503433334
$newValues[2,13] = 0.4943960553090001
$newValues[2,14] = 0.00005081536077573487
$newValues[2,15] = 0.00005081536077573488

# Row 5
$newValues[3,0] = 3
$newValues[3,1] = 1
$newValues[3,2] = 1.208813
$newValues[3,3] = 3.626439
$newValues[3,4] = 0.06175614228936263
$newValues[3,5] = 0.06175614228936264
$newValues[3,6] = 3
$newValues[3,7] = 1
$newValues[3,8] = 54.908578
$newValues[3,9] = 164.725734
$newValues[3,10] = 0.9942181455359038
$newValues[3,11] = 0.9942181455359038
$newValues[3,12] = 66.374202897914
$newValues[3,13] = 597.367826081226
$newValues[3,14] = 0.06139907726238152
$newValues[3,15] = 0.06139907726238153

# Row 6
$newValues[4,0] = 3
$newValues[4,1] = 1
$newValues[4,2] = 2.407869333333334
$newValues[4,3] = 7.223608
$newValues[4,4] = 0.1230138335404451
$newValues[4,5] = 0.1230138335404451
$newValues[4,6] = 3
$newValues[4,7] = 1
$newValues[4,8] = 0.1956483333333333
$newValues[4,9] = 0.5869450000000001
$newValues[4,10] = 0.003542563479678112
$newValues[4,11] = 0.003542563479678113
$newValues[4,12] = 0.4710956219511112
$newValues[4,13] = 4.239860597560001
$newValues[4,14] = 0.0004357843141955833
$newValues[4,15] = 0.0004357843141955833

# Row 7
$newValues[5,0] = 3
$newValues[5,1] = 1
$newValues[5,2] = 2.407869333333334
$newValues[5,3] = 7.223608
$newValues[5,4] = 0.1230138335404451
$newValues[5,5] = 0.1230138335404451
$newValues[5,6] = 1
$newValues[5,7] = 0.3333333333333333
$newValues[5,8] = 0.07822766666666667
$newValues[5,9] = 0.234683
$newValues[5,10] = 0.001416452010156485
$newValues[5,11] = 0.001416452010156486
$newValues[5,12] = 0.1883619995848889
$newValues[5,13] = 1.695257996264
$newValues[5,14] = 0.0001742431917954187
$newValues[5,15] = 0.0001742431917954188

# Row 8
$newValues[6,0] = 3
$newValues[6,1] = 1
$newValues[6,2] = 2.407869333333334
$newValues[6,3] = 7.223608
$newValues[6,4] = 0.1230138335404451
$newValues[6,5] = 0.1230138335404451
$newValues[6,6] = 2
$newValues[6,7] = 0.6666666666666666
$newValues[6,8] = 0.04544366666666667
$newValues[6,9] = 0.136331
$newValues[6,10] = 0.0008228389742616373
$newValues[6,11] = 0.0008228389742616374
$newValues[6,12] = 0.1094224113608889
$newValues[6,13] = 0.9848017022480001
$newValues[6,14] = 0.0001012205766104116
$newValues[6,15] = 0.0001012205766104117

# Row 9
$newValues[7,0] = 3
$newValues[7,1] = 1
$newValues[7,2] = 2.407869333333334
$newValues[7,3] = 7.223608
$newValues[7,4] = 0.1230138335404451
$newValues[7,5] = 0.1230138335404451
$newValues[7,6] = 3
$newValues[7,7] = 1
$newValues[7,8] = 54.908578
$newValues[7,9] = 164.725734
$newValues[7,10] = 0.9942181455359038
$newValues[7,11] = 0.9942181455359038
$newValues[7,12] = 132.2126811031413
$newValues[7,13] = 1189.914129928272
$newValues[7,14] = 0.1223025854578437
$newValues[7,15] = 0.1223025854578437

# Row 10
$newValues[8,0] = 3
$newValues[8,1] = 1
$newValues[8,2] = 0.9221053333333332
$newValues[8,3] = 2.766316
$newValues[8,4] = 0.04710874897201923
$newValues[8,5] = 0.04710874897201923
$newValues[8,6] = 3
$newValues[8,7] = 1
$newValues[8,8] = 0.1956483333333333
$newValues[8,9] = 0.5869450000000001
$newValues[8,10] = 0.003542563479678112
$newValues[8,11] = 0.003542563479678113
$newValues[8,12] = 0.1804083716244444
$newValues[8,13] = 1.62367534462
$newValues[8,14] = 0.0001668857336815991
$newValues[8,15] = 0.0001668857336815992

# Row 11
$newValues[9,0] = 3
$newValues[9,1] = 1
$newValues[9,2] = 0.9221053333333332
$newValues[9,3] = 2.766316
$newValues[9,4] = 0.04710874897201923
$newValues[9,5] = 0.04710874897201923
$newValues[9,6] = 1
$newValues[9,7] = 0.3333333333333333
$newValues[9,8] = 0.07822766666666667
$newValues[9,9] = 0.234683
$newValues[9,10] = 0.001416452010156485
$newValues[9,11] = 0.001416452010156486
$newValues[9,12] = 0.07213414864755555
$newValues[9,13] = 0.649207337828
$newValues[9,14] = 0.0000667272821773739
$newValues[9,15] = 0.00006672728217737391

# Row 12
$newValues[10,0] = 3
$newValues[10,1] = 1
$newValues[10,2] = 0.9221053333333332
$newValues[10,3] = 2.766316
$newValues[10,4] = 0.04710874897201923
$newValues[10,5] = 0.04710874897201923
$newValues[10,6] = 2
$newValues[10,7] = 0.6666666666666666
$newValues[10,8] = 0.04544366666666667
$newValues[10,9] = 0.136331
$newValues[10,10] = 0.0008228389742616373
$newValues[10,11] = 0.0008228389742616374
$newValues[10,12] = 0.04190384739955555
$newValues[10,13] = 0.377134626596
$newValues[10,14] = 0.00003876291468288526
$newValues[10,15] = 0.00003876291468288527

# Row 13
$newValues[11,0] = 3
$newValues[11,1] = 1
$newValues[11,2] = 0.9221053333333332
$newValues[11,3] = 2.766316
$newValues[11,4] = 0.04710874897201923
$newValues[11,5] = 0.04710874897201923
$newValues[11,6] = 3
$newValues[11,7] = 1
$newValues[11,8] = 54.908578
$newValues[11,9] = 164.725734
$newValues[11,10] = 0.9942181455359038
$newValues[11,11] = 0.9942181455359038
$newValues[11,12] = 50.63149261954933
$newValues[11,13] = 455.6834335759439
$newValues[11,14] = 0.04683637304147737
$newValues[11,15] = 0.04683637304147738

# Row 14
$newValues[12,0] = 3
$newValues[12,1] = 1
$newValues[12,2] = 15.03518433333333
$newValues[12,3] = 45.105553
$newValues[12,4] = 0.768121275198173
$newValues[12,5] = 0.7681212751981731
$newValues[12,6] = 3
$newValues[12,7] = 1
$newValues[12,8] = 0.1956483333333333
$newValues[12,9] = 0.5869450000000001
$newValues[12,10] = 0.003542563479678112
$newValues[12,11] = 0.003542563479678113
$newValues[12,12] = 2.941608756176111
$newValues[12,13] = 26.474478805585
$newValues[12,14] = 0.002721118377480829
$newValues[12,15] = 0.002721118377480829

# Row 15
$newValues[13,0] = 3
$newValues[13,1] = 1
$newValues[13,2] = 15.03518433333333
$newValues[13,3] = 45.105553
$newValues[13,4] = 0.768121275198173
$newValues[13,5] = 0.7681212751981731
$newValues[13,6] = 1
$newValues[13,7] = 0.3333333333333333
$newValues[13,8] = 0.07822766666666667
$newValues[13,9] = 0.234683
$newValues[13,10] = 0.001416452010156485
$newValues[13,11] = 0.001416452010156486
$newValues[13,12] = 1.176167388299889
$newValues[13,13] = 10.585506494699
$newValues[13,14] = 0.001088006924298415
$newValues[13,15] = 0.001088006924298415

# Row 16
$newValues[14,0] = 3
$newValues[14,1] = 1
$newValues[14,2] = 15.03518433333333
$newValues[14,3] = 45.105553
$newValues[14,4] = 0.768121275198173
$newValues[14,5] = 0.7681212751981731
$newValues[14,6] = 2
$newValues[14,7] = 0.6666666666666666
$newValues[14,8] = 0.04544366666666667
$newValues[14,9] = 0.136331
$newValues[14,10] = 0.0008228389742616373
$newValues[14,11] = 0.0008228389742616374
$newValues[14,12] = 0.6832539051158889
$newValues[14,13] = 6.149285146043001
$newValues[14,14] = 0.0006320401221926055
$newValues[14,15] = 0.0006320401221926056

# Row 17
$newValues[15,0] = 3
$newValues[15,1] = 1
$newValues[15,2] = 15.03518433333333
$newValues[15,3] = 45.105553
$newValues[15,4] = 0.768121275198173
$newValues[15,5] = 0.7681212751981731
$newValues[15,6] = 3
$newValues[15,7] = 1
$newValues[15,8] = 54.908578
$newValues[15,9] = 164.725734
$newValues[15,10] = 0.9942181455359038
$newValues[15,11] = 0.9942181455359038
$newValues[15,12] = 825.5605917112114
$newValues[15,13] = 7430.045325400902
$newValues[15,14] = 0.7636801097742012
$newValues[15,15] = 0.7636801097742013

$ws.Range("E2:T17").Value2 = $newValues

